$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The four small "UndoStack/RedoStack" tables on the diagram each contain a
# row referencing the old "prevAddressBook" variable name. Rename it to
# "prevImdb" in each, leaving the rest of the cell text (" = s1"/" = s3")
# untouched.
#
# Shape ids (per p:cNvPr) -> "Table 58" (id 59), "Table 62" (id 63),
# "Table 23" (id 24), "Table 24" (id 25). These are shapes 6-9 on the slide.
$tableShapeIndexes = @(6, 7, 8, 9)

foreach ($idx in $tableShapeIndexes) {
    $shape = $s.Shapes.Item($idx)
    $tbl = $shape.Table

    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        $cellRange = $tbl.Cell($r, 1).Shape.TextFrame.TextRange
        $runs = $cellRange.Runs()
        for ($j = 1; $j -le $runs.Count; $j++) {
            $run = $runs.Item($j)
            if ($run.Text -like "prevAddressBook*") {
                $run.Text = $run.Text -replace "prevAddressBook", "prevImdb"
            }
        }
    }
}
